$wb = $excel.ActiveWorkbook

# --- Preserve/establish selection on "Original features" sheet (G28) before switching active sheet ---
$wb.Worksheets.Item("Original features").Range("G28").Select()

# --- Create the new worksheet as the last tab, forcing sheetId=6 (matches diff) ---
# NOTE: worksheet handles returned by Add() are resolved by POSITION, not stable identity, so we
# must rename each newly-added sheet immediately (before the next structural change shifts indices).
# The engine assigns sheetId = max(existing sheetIds)+1 per Add(), and frees the id back to the pool
# on Delete(); adding a throwaway sheet first (claiming id 5), then the real one (claiming id 6), then
# deleting the throwaway leaves the real sheet holding sheetId 6.
$throwaway = $wb.Worksheets.Add()
$throwaway.Name = "zzz_throwaway1"

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "new features importance2"

$wb.Worksheets.Item("zzz_throwaway1").Delete()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wb.Worksheets.Item("new features importance2").Move([System.Reflection.Missing]::Value, $lastSheet)
$wb.Worksheets.Item("new features importance2").Activate()

# --- From here on, no further structural (add/delete/move) changes happen, so it's safe to hold a reference ---
$ws = $wb.Worksheets.Item("new features importance2")

# --- Column widths matching the diff's <cols> ---
$ws.Columns("A:A").ColumnWidth = 18.1640625
$ws.Columns("B:C").ColumnWidth = 12.1640625

# --- Styling: borders across A1:C42, number format on B:C, black font on header A1 ---
$ws.Range("A1:C42").Borders.LineStyle = 1
$ws.Range("B1:C42").NumberFormat = "0.000"
$ws.Range("A1").Font.Color = 0

$ws.Range('A1').Value = 'LightGBM_Features'
$ws.Range('B1').Value = 'Gain'
$ws.Range('C1').Value = 'Frequency'
$ws.Range('E1').Value = 'ps_car_10_catnew'
$ws.Range('A2').Value = 'ps_ind_05_catnew'
$ws.Range('B2').Value = 0.101656934296827
$ws.Range('C2').Value = 0.066965699877092305
$ws.Range('E2').Value = 'ps_ind_14'
$ws.Range('A3').Value = 'avg_car13_on_ind02'
$ws.Range('B3').Value = 0.096344092634208303
$ws.Range('C3').Value = 0.035123731326263001
$ws.Range('A4').Value = 'ps_ind_03'
$ws.Range('B4').Value = 0.088698965442326402
$ws.Range('C4').Value = 0.087695290226935799
$ws.Range('A5').Value = 'ps_car_11_catnew'
$ws.Range('B5').Value = 0.084505237294001001
$ws.Range('C5').Value = 0.046274122223489297
$ws.Range('A6').Value = 'ps_reg_03'
$ws.Range('B6').Value = 0.067088878201772104
$ws.Range('C6').Value = 0.072806984199389302
$ws.Range('A7').Value = 'ps_ind_17_bin'
$ws.Range('B7').Value = 0.062213766387722899
$ws.Range('C7').Value = 0.0301694099162454
$ws.Range('A8').Value = 'ps_ind_15'
$ws.Range('B8').Value = 0.045318318014859603
$ws.Range('C8').Value = 0.056892335373348001
$ws.Range('A9').Value = 'ps_car_01_catnew'
$ws.Range('B9').Value = 0.039305506223947399
$ws.Range('C9').Value = 0.033818628755337601
$ws.Range('A10').Value = 'ps_car_03_catnew'
$ws.Range('B10').Value = 0.033232462402815299
$ws.Range('C10').Value = 0.024315454695201501
$ws.Range('A11').Value = 'ps_reg_01'
$ws.Range('B11').Value = 0.032483758386462801
$ws.Range('C11').Value = 0.042713599675624998
$ws.Range('A12').Value = 'ps_ind_06_bin'
$ws.Range('B12').Value = 0.029993466073982299
$ws.Range('C12').Value = 0.0131523928992283
$ws.Range('A13').Value = 'ps_ind_01'
$ws.Range('B13').Value = 0.027603053887888002
$ws.Range('C13').Value = 0.0473131359207309
$ws.Range('A14').Value = 'ps_reg_02'
$ws.Range('B14').Value = 0.0252971077138569
$ws.Range('C14').Value = 0.033438501792932201
$ws.Range('A15').Value = 'ps_car_07_catnew'
$ws.Range('B15').Value = 0.023395824042367799
$ws.Range('C15').Value = 0.021325122590945399
$ws.Range('A16').Value = 'ps_car_09_catnew'
$ws.Range('B16').Value = 0.021389975640313301
$ws.Range('C16').Value = 0.0229469976305419
$ws.Range('A17').Value = 'avg_car13_on_car11'
$ws.Range('B17').Value = 0.0197792005363027
$ws.Range('C17').Value = 0.043017701245549399
$ws.Range('A18').Value = 'ps_ind_16_bin'
$ws.Range('B18').Value = 0.018332561813045101
$ws.Range('C18').Value = 0.0100987063012379
$ws.Range('A19').Value = 'ps_ind_07_bin'
$ws.Range('B19').Value = 0.0182159685046967
$ws.Range('C19').Value = 0.0103267824786812
$ws.Range('A20').Value = 'avg_car13_on_ind04'
$ws.Range('B20').Value = 0.0171752669295837
$ws.Range('C20').Value = 0.014482837267647399
$ws.Range('A21').Value = 'avg_car13_on_ind05'
$ws.Range('B21').Value = 0.015341075746750499
$ws.Range('C21').Value = 0.012645556949354399
$ws.Range('A22').Value = 'ps_ind_02_catnew'
$ws.Range('B22').Value = 0.013781024809619799
$ws.Range('C22').Value = 0.018600879360373
$ws.Range('A23').Value = 'avg_car13_on_car04'
$ws.Range('B23').Value = 0.011465591000674199
$ws.Range('C23').Value = 0.025835962544823302
$ws.Range('A24').Value = 'ps_ind_09_bin'
$ws.Range('B24').Value = 0.0113724942097926
$ws.Range('C24').Value = 0.010820947529808301
$ws.Range('A25').Value = 'ps_car_15'
$ws.Range('B25').Value = 0.011337150692633
$ws.Range('C25').Value = 0.017207080498219701
$ws.Range('A26').Value = 'avg_car13_on_car01'
$ws.Range('B26').Value = 0.011006266733819801
$ws.Range('C26').Value = 0.026013355127279199
$ws.Range('A27').Value = 'car13_car15'
$ws.Range('B27').Value = 0.0096431304970415292
$ws.Range('C27').Value = 0.0263427984946972
$ws.Range('A28').Value = 'ps_car_14'
$ws.Range('B28').Value = 0.0093597726642553596
$ws.Range('C28').Value = 0.027787280951837901
$ws.Range('A29').Value = 'avg_car13_on_car02'
$ws.Range('B29').Value = 0.0079618701715402796
$ws.Range('C29').Value = 0.021502515173401199
$ws.Range('A30').Value = 'avg_car13_on_car09'
$ws.Range('B30').Value = 0.0066362163637294497
$ws.Range('C30').Value = 0.018765601044082099
$ws.Range('A31').Value = 'ps_car_13'
$ws.Range('B31').Value = 0.0066296047960354597
$ws.Range('C31').Value = 0.0080586916029953999
$ws.Range('A32').Value = 'ps_car_11'
$ws.Range('B32').Value = 0.0064159892816749702
$ws.Range('C32').Value = 0.016180737699725001
$ws.Range('A33').Value = 'ps_ind_04_catnew'
$ws.Range('B33').Value = 0.0059304580792587302
$ws.Range('C33').Value = 0.0094524904651486896
$ws.Range('A34').Value = 'ps_ind_08_bin'
$ws.Range('B34').Value = 0.0058070424355452597
$ws.Range('C34').Value = 0.0075518556531214802
$ws.Range('A35').Value = 'avg_car13_on_car07'
$ws.Range('B35').Value = 0.0053690981965173398
$ws.Range('C35').Value = 0.013431152671658999
$ws.Range('A36').Value = 'ps_car_06_catnew'
$ws.Range('B36').Value = 0.0036175072691710698
$ws.Range('C36').Value = 0.0102127443899596
$ws.Range('A37').Value = 'ps_car_12'
$ws.Range('B37').Value = 0.0019417804630933999
$ws.Range('C37').Value = 0.0060947022972339402
$ws.Range('A38').Value = 'ps_car_05_catnew'
$ws.Range('B38').Value = 0.0012060050678853701
$ws.Range('C38').Value = 0.0032944336741805099
$ws.Range('A39').Value = 'ps_car_08_catnew'
$ws.Range('B39').Value = 0.00117362041850181
$ws.Range('C39').Value = 0.0025215088506227701
$ws.Range('A40').Value = 'ps_ind_18_bin'
$ws.Range('B40').Value = 0.0010738444852093001
$ws.Range('C40').Value = 0.0021540527869641801
$ws.Range('A41').Value = 'ps_car_04_catnew'
$ws.Range('B41').Value = 0.00076897910828771499
$ws.Range('C41').Value = 0.0021033691919767898
$ws.Range('A42').Value = 'ps_car_02_catnew'
$ws.Range('B42').Value = 0.00013113308198388899
$ws.Range('C42').Value = 0.00054484864611446901
